$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (accuracy) values for existing rows 2-116, and extend with new rows 117-118
$bValues = @{
    2 = 0.796875
    3 = 0.71875
    4 = 0.65625
    5 = 0.59375
    6 = 0.53125
    7 = 0.4375
    8 = 0.453125
    9 = 0.421875
    10 = 0.46875
    11 = 0.46875
    12 = 0.484375
    13 = 0.5
    14 = 0.484375
    15 = 0.5
    16 = 0.5
    17 = 0.5
    18 = 0.5
    19 = 0.5
    20 = 0.5
    21 = 0.5
    22 = 0.5
    23 = 0.515625
    24 = 0.515625
    25 = 0.515625
    26 = 0.515625
    27 = 0.5
    28 = 0.5
    29 = 0.484375
    30 = 0.5
    31 = 0.5
    32 = 0.5
    33 = 0.484375
    34 = 0.484375
    35 = 0.5
    36 = 0.46875
    37 = 0.46875
    38 = 0.46875
    39 = 0.46875
    40 = 0.46875
    41 = 0.46875
    42 = 0.46875
    43 = 0.46875
    44 = 0.484375
    45 = 0.46875
    46 = 0.484375
    47 = 0.484375
    48 = 0.484375
    49 = 0.484375
    50 = 0.484375
    51 = 0.484375
    52 = 0.484375
    53 = 0.46875
    54 = 0.484375
    55 = 0.46875
    56 = 0.484375
    57 = 0.484375
    58 = 0.484375
    59 = 0.484375
    60 = 0.46875
    61 = 0.484375
    62 = 0.484375
    63 = 0.484375
    64 = 0.484375
    65 = 0.484375
    66 = 0.484375
    67 = 0.484375
    68 = 0.484375
    69 = 0.484375
    70 = 0.484375
    71 = 0.484375
    72 = 0.484375
    73 = 0.484375
    74 = 0.484375
    75 = 0.484375
    76 = 0.484375
    77 = 0.484375
    78 = 0.484375
    79 = 0.484375
    80 = 0.484375
    81 = 0.484375
    82 = 0.484375
    83 = 0.484375
    84 = 0.484375
    85 = 0.484375
    86 = 0.484375
    87 = 0.484375
    88 = 0.484375
    89 = 0.484375
    90 = 0.484375
    91 = 0.484375
    92 = 0.484375
    93 = 0.484375
    94 = 0.484375
    95 = 0.484375
    96 = 0.484375
    97 = 0.484375
    98 = 0.484375
    99 = 0.484375
    100 = 0.484375
    101 = 0.484375
    102 = 0.484375
    104 = 0.390625
    105 = 0.546875
    106 = 0.34375
    107 = 0.359375
    108 = 0.421875
    109 = 0.484375
    110 = 0.4375
    111 = 0.46875
    112 = 0.671875
    113 = 0.46875
    114 = 0.328125
    115 = 0.46875
    116 = 0.375
    117 = 0.4375
    118 = 0.360655737704918
}

foreach ($row in $bValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $bValues[$row]
}

# New rows 117 and 118: column A holds the same repr string as the other DisplayOutputs rows
$displayOutputsRepr = "<__main__.DisplayOutputs object at 0x7f16c0289cd0>"
$ws.Cells.Item(117, 1).Value = $displayOutputsRepr
$ws.Cells.Item(118, 1).Value = $displayOutputsRepr

# Update the existing DisplayOutputs repr strings in column A (rows 102-116) to the new memory address
for ($r = 102; $r -le 116; $r++) {
    $ws.Cells.Item($r, 1).Value = $displayOutputsRepr
}

# Restore the active cell / selection to match the authored view state
$ws.Range("A2:B116").Select()
